$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title / header text updates (Volume number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Weekly crime statistics grid updates (rows 14-30) ---
# Row 14
$ws.Range("C14").Value = 1
$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 5
$ws.Range("H14").Value = 25
$ws.Range("I14").Value = 4
$ws.Range("J14").Value = 3
$ws.Range("K14").Value = 33.333333333333
$ws.Range("L14").Value = -42.857142857142
$ws.Range("M14").Value = -20
$ws.Range("N14").Value = -87.878787878787

# Row 15
$ws.Range("D15").Value = 7
$ws.Range("E15").Value = -14.285714285714
$ws.Range("F15").Value = 27
$ws.Range("G15").Value = 16
$ws.Range("H15").Value = 68.75
$ws.Range("I15").Value = 22
$ws.Range("J15").Value = 13
$ws.Range("K15").Value = 69.230769230769
$ws.Range("L15").Value = 46.666666666666
$ws.Range("M15").Value = 100

# Row 16
$ws.Range("C16").Value = 36
$ws.Range("E16").Value = -5.263157894736
$ws.Range("F16").Value = 157
$ws.Range("G16").Value = 182
$ws.Range("H16").Value = -13.736263736263
$ws.Range("I16").Value = 133
$ws.Range("J16").Value = 130
$ws.Range("K16").Value = 2.307692307692
$ws.Range("L16").Value = 31.683168316831
$ws.Range("M16").Value = -40.090090090090
$ws.Range("N16").Value = -87.333333333333

# Row 17
$ws.Range("D17").Value = 75
$ws.Range("E17").Value = -16
$ws.Range("F17").Value = 265
$ws.Range("G17").Value = 261
$ws.Range("H17").Value = 1.532567049808
$ws.Range("I17").Value = 216
$ws.Range("J17").Value = 211
$ws.Range("K17").Value = 2.369668246445
$ws.Range("L17").Value = 16.129032258064
$ws.Range("M17").Value = 24.855491329479
$ws.Range("N17").Value = -54.811715481171

# Row 18
$ws.Range("C18").Value = 47
$ws.Range("D18").Value = 45
$ws.Range("E18").Value = 4.444444444444
$ws.Range("F18").Value = 209
$ws.Range("G18").Value = 223
$ws.Range("H18").Value = -6.278026905829
$ws.Range("I18").Value = 160
$ws.Range("J18").Value = 160
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 20.300751879699
$ws.Range("M18").Value = -4.191616766467
$ws.Range("N18").Value = -76.505139500734

# Row 19
$ws.Range("C19").Value = 124
$ws.Range("D19").Value = 106
$ws.Range("E19").Value = 16.981132075471
$ws.Range("F19").Value = 433
$ws.Range("G19").Value = 387
$ws.Range("H19").Value = 11.886304909560
$ws.Range("I19").Value = 337
$ws.Range("J19").Value = 310
$ws.Range("K19").Value = 8.709677419354
$ws.Range("L19").Value = 27.169811320754
$ws.Range("M19").Value = 77.368421052631
$ws.Range("N19").Value = -8.423913043478

# Row 20
$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 33
$ws.Range("E20").Value = -12.121212121212
$ws.Range("F20").Value = 123
$ws.Range("G20").Value = 139
$ws.Range("H20").Value = -11.510791366906
$ws.Range("I20").Value = 97
$ws.Range("J20").Value = 115
$ws.Range("K20").Value = -15.652173913043
$ws.Range("L20").Value = 46.969696969697
$ws.Range("M20").Value = 22.784810126582
$ws.Range("N20").Value = -84.329563812601

# Row 21
$ws.Range("C21").Value = 306
$ws.Range("D21").Value = 305
$ws.Range("E21").Value = 0.327868852459
$ws.Range("F21").Value = 1219
$ws.Range("G21").Value = 1212
$ws.Range("H21").Value = 0.577557755775
$ws.Range("I21").Value = 969
$ws.Range("J21").Value = 942
$ws.Range("K21").Value = 2.866242038216
$ws.Range("L21").Value = 25.355756791720
$ws.Range("M21").Value = 14.403778040141
$ws.Range("N21").Value = -70.294297976701

# Row 22
$ws.Range("C22").Value = 5
$ws.Range("D22").Value = 11
$ws.Range("E22").Value = -54.545454545454
$ws.Range("F22").Value = 25
$ws.Range("G22").Value = 25
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 20
$ws.Range("J22").Value = 22
$ws.Range("K22").Value = -9.090909090909
$ws.Range("L22").Value = 81.818181818181
$ws.Range("M22").Value = -4.761904761904

# Row 23
$ws.Range("C23").Value = 30
$ws.Range("E23").Value = 3.448275862068
$ws.Range("F23").Value = 116
$ws.Range("G23").Value = 128
$ws.Range("H23").Value = -9.375
$ws.Range("I23").Value = 94
$ws.Range("J23").Value = 92
$ws.Range("K23").Value = 2.173913043478
$ws.Range("L23").Value = 28.767123287671
$ws.Range("M23").Value = 80.769230769230

# Row 24
$ws.Range("C24").Value = 224
$ws.Range("D24").Value = 224
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 888
$ws.Range("G24").Value = 853
$ws.Range("H24").Value = 4.103165298944
$ws.Range("I24").Value = 685
$ws.Range("J24").Value = 630
$ws.Range("K24").Value = 8.730158730158
$ws.Range("L24").Value = 18.512110726643
$ws.Range("M24").Value = 15.709459459459

# Row 25
$ws.Range("C25").Value = 117
$ws.Range("D25").Value = 108
$ws.Range("E25").Value = 8.333333333333
$ws.Range("F25").Value = 447
$ws.Range("G25").Value = 400
$ws.Range("H25").Value = 11.75
$ws.Range("I25").Value = 347
$ws.Range("J25").Value = 314
$ws.Range("K25").Value = 10.509554140127
$ws.Range("L25").Value = 39.919354838709
$ws.Range("M25").Value = -20.776255707762

# Row 26
$ws.Range("C26").Value = 7
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -12.5
$ws.Range("F26").Value = 33
$ws.Range("G26").Value = 25
$ws.Range("H26").Value = 32
$ws.Range("I26").Value = 27
$ws.Range("J26").Value = 22
$ws.Range("K26").Value = 22.727272727272
$ws.Range("L26").Value = 3.846153846153

# Row 27
$ws.Range("C27").Value = 12
$ws.Range("D27").Value = 5
$ws.Range("E27").Value = 140
$ws.Range("F27").Value = 47
$ws.Range("G27").Value = 33
$ws.Range("H27").Value = 42.424242424242
$ws.Range("I27").Value = 41
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = 70.833333333333
$ws.Range("L27").Value = 51.851851851851

# Row 28
$ws.Range("C28").Value = 2
$ws.Range("D28").Value = 4
$ws.Range("E28").Value = -50
$ws.Range("F28").Value = 9
$ws.Range("G28").Value = 14
$ws.Range("H28").Value = -35.714285714285
$ws.Range("I28").Value = 8
$ws.Range("J28").Value = 12
$ws.Range("K28").Value = -33.333333333333
$ws.Range("L28").Value = -33.333333333333
$ws.Range("M28").Value = -61.904761904761
$ws.Range("N28").Value = -93.495934959349

# Row 29
$ws.Range("C29").Value = 2
$ws.Range("D29").Value = 4
$ws.Range("E29").Value = -50
$ws.Range("F29").Value = 9
$ws.Range("G29").Value = 13
$ws.Range("H29").Value = -30.769230769230
$ws.Range("I29").Value = 8
$ws.Range("J29").Value = 11
$ws.Range("K29").Value = -27.272727272727
$ws.Range("L29").Value = -33.333333333333
$ws.Range("M29").Value = -52.941176470588
$ws.Range("N29").Value = -92.982456140350

# Row 30
$ws.Range("D30").Value = 3
$ws.Range("E30").Value = -100
$ws.Range("F30").Value = 7
$ws.Range("G30").Value = 5
$ws.Range("H30").Value = 40
$ws.Range("I30").Value = 6
$ws.Range("J30").Value = 5
$ws.Range("K30").Value = 20
$ws.Range("L30").Value = 100

# D30 and E30 change from text placeholders to numeric cells; set number formats to match target styles
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
